$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "G3"
$ws.Range("B4").Value = "Test2"
$ws.Range("C4").Value = "Daily"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 45860
$ws.Range("E4").Style = $ws.Range("E3").Style
$ws.Range("E4").NumberFormat = $ws.Range("E3").NumberFormat
$ws.Range("F4").Value = 30
